$d = $word.ActiveDocument

# Update the date heading paragraph.
$d.Content.Find.Execute("2026-01-21 Wednesday", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "2026-01-22 Thursday", 2)

# Update the practice-problem table cells. Cells are addressed by
# (row, column) to avoid ambiguity between duplicate problem texts
# (e.g. "13÷4=" appears twice with two different replacements).
$tbl = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "84÷4=" },
    @{ Row = 1;  Col = 2; New = "54÷4=" },
    @{ Row = 1;  Col = 3; New = "51÷7=" },
    @{ Row = 1;  Col = 4; New = "38÷2=" },
    @{ Row = 1;  Col = 5; New = "99÷4=" },

    @{ Row = 5;  Col = 1; New = "87÷8=" },
    @{ Row = 5;  Col = 2; New = "85÷6=" },
    @{ Row = 5;  Col = 3; New = "22÷8=" },
    @{ Row = 5;  Col = 4; New = "50÷3=" },
    @{ Row = 5;  Col = 5; New = "47÷9=" },

    @{ Row = 9;  Col = 1; New = "19÷9=" },
    @{ Row = 9;  Col = 2; New = "92÷6=" },
    @{ Row = 9;  Col = 3; New = "13÷7=" },
    @{ Row = 9;  Col = 4; New = "59÷3=" },
    @{ Row = 9;  Col = 5; New = "44÷9=" },

    @{ Row = 13; Col = 1; New = "38÷3=" },
    @{ Row = 13; Col = 2; New = "10÷6=" },
    @{ Row = 13; Col = 3; New = "36÷5=" },
    @{ Row = 13; Col = 4; New = "29÷3=" },
    @{ Row = 13; Col = 5; New = "89÷3=" },

    @{ Row = 17; Col = 1; New = "74÷8=" },
    @{ Row = 17; Col = 2; New = "63÷7=" },
    @{ Row = 17; Col = 3; New = "64÷3=" },
    @{ Row = 17; Col = 4; New = "99÷8=" },
    @{ Row = 17; Col = 5; New = "77÷8=" }
)

foreach ($u in $updates) {
    $cell = $tbl.Rows.Item($u.Row).Cells.Item($u.Col)
    $cell.Range.Text = $u.New
}

Write-Host "done"
